$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("Sprint3")
$r = $ws7.Range("Z2")
$r.Value = "x"
$r.BorderAround(1)
